$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "69.304.97"
Set-TextValue "E2" "  +0.82%  "
Set-TextValue "D3" "3.788.11"
Set-TextValue "E3" "  +1.07%  "
Set-TextValue "E4" "  +0.16%  "
Set-TextValue "D5" "603.32"
Set-TextValue "E5" "  +0.26%  "
Set-TextValue "D6" "165.01"
Set-TextValue "E6" "  -2.41%  "
Set-TextValue "D7" "3.785.83"
Set-TextValue "E7" "  +1.03%  "
Set-TextValue "E8" "  -0.01%  "
Set-TextValue "E9" "  +0.63%  "
Set-TextValue "E10" "  +4.25%  "
Set-TextValue "D11" "6.32"
Set-TextValue "E11" "  -0.22%  "
Set-TextValue "E12" "  -0.41%  "
Set-TextValue "D13" "37.60"
Set-TextValue "E13" "  -2.04%  "
Set-TextValue "E14" "  -0.40%  "
Set-TextValue "D15" "4.421.81"
Set-TextValue "E15" "  +1.07%  "
Set-TextValue "D16" "3.786.56"
Set-TextValue "E16" "  +0.97%  "
Set-TextValue "D17" "69.405.87"
Set-TextValue "E17" "  +0.94%  "
Set-TextValue "D18" "7.43"
Set-TextValue "E18" "  +1.85%  "
Set-TextValue "D19" "17.62"
Set-TextValue "E19" "  +3.12%  "
Set-TextValue "E20" "  -0.98%  "
Set-TextValue "D21" "11.36"
Set-TextValue "E21" "  +4.90%  "
Set-TextValue "D22" "492.97"
Set-TextValue "E22" "  -0.70%  "
Set-TextValue "D23" "0.725"
Set-TextValue "E23" "  -0.57%  "
Set-TextValue "D24" "0.0000151"
Set-TextValue "E24" "  -2.91%  "
Set-TextValue "D25" "84.89"
Set-TextValue "E25" "  -0.61%  "
Set-TextValue "E26" "  -2.71%  "
Set-TextValue "D27" "12.31"
Set-TextValue "E27" "  -0.23%  "
Set-TextValue "D28" "10.12"
Set-TextValue "E28" "  -2.06%  "
Set-TextValue "E29" "  +0.07%  "
Set-TextValue "D30" "2.98"
Set-TextValue "E30" "  -0.35%  "
Set-TextValue "D31" "8.11"
Set-TextValue "E31" "  +2.10%  "
Set-TextValue "E32" "  -3.77%  "
Set-TextValue "D33" "31.92"
Set-TextValue "E33" "  +0.10%  "
Set-TextValue "D34" "3.932.28"
Set-TextValue "E34" "  +1.01%  "
Set-TextValue "D35" "3.734.50"
Set-TextValue "E35" "  +1.40%  "
Set-TextValue "E36" "  -0.77%  "
Set-TextValue "D37" "5.97"
Set-TextValue "E37" "  +1.80%  "
Set-TextValue "B38" "Kaspa"
Set-TextValue "C38" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D38" "0.140"
Set-TextValue "E38" "  +4.89%  "
Set-TextValue "B39" "Mantle"
Set-TextValue "C39" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D39" "1.02"
Set-TextValue "E39" "  +0.23%  "
Set-TextValue "E40" "  +0.08%  "
Set-TextValue "D41" "0.325"
Set-TextValue "E41" "  +0.32%  "
Set-TextValue "D42" "3.06"
Set-TextValue "E42" "  +4.24%  "
Set-TextValue "B43" "OKB"
Set-TextValue "C43" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D43" "48.50"
Set-TextValue "E43" "  -0.56%  "
Set-TextValue "B44" "Stacks"
Set-TextValue "C44" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D44" "1.99"
Set-TextValue "E44" "  +0.72%  "
Set-TextValue "D45" "423.98"
Set-TextValue "E45" "  -3.37%  "
Set-TextValue "D46" "8.44"
Set-TextValue "E46" "  -0.40%  "
Set-TextValue "D48" "40.34"
Set-TextValue "E48" "  -0.64%  "
Set-TextValue "D49" "141.84"
Set-TextValue "E49" "  +0.21%  "
Set-TextValue "D50" "2.814.43"
Set-TextValue "E50" "  +0.84%  "
Set-TextValue "E51" "  +6.94%  "
